$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 720.8125
$ws.Range("I28").Value = 668.86664
$ws.Range("K28").Value = 668.86664
$ws.Range("M28").Value = -183.86664

$ws.Range("H33").Value = 190.76923
$ws.Range("I33").Value = 190.76923
$ws.Range("K33").Value = 190.76923
$ws.Range("M33").Value = 38.23077000000001

$ws.Range("H76").Value = 1000
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -685
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 1000
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 92
$ws.Range("N79").ClearContents()

$ws.Range("H132").Value = 5931.625
$ws.Range("I132").Value = 3491.2
$ws.Range("K132").Value = 10473.6
$ws.Range("M132").Value = -7943.599999999999

$ws.Range("H138").Value = 8657.842000000001
$ws.Range("I138").Value = 2499.6667
$ws.Range("J138").Value = 9812.5
$ws.Range("K138").Value = 7499.000100000001
$ws.Range("L138").Value = 29437.5
$ws.Range("M138").Value = -2359.000100000001
$ws.Range("N138").Value = -39717.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 57.5
$ws.Range("I5").Value = 43.333332
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 43.333332
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 68.666668
$ws.Range("N5").Value = -324

$ws.Range("H32").Value = 5195.706
$ws.Range("I32").Value = 4809.4287
$ws.Range("K32").Value = 4809.4287
$ws.Range("M32").Value = -4522.4287

$ws.Range("H45").Value = 2823
$ws.Range("I45").Value = 1734.5
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 1734.5
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -1357.5
$ws.Range("N45").Value = -5754

$ws.Range("H110").Value = 3500
$ws.Range("I110").Value = 3500
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 3500
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1455
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 2000.25
$ws.Range("I122").Value = 2000.25
$ws.Range("K122").Value = 6000.75
$ws.Range("M122").Value = -3550.75

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 57.5
$ws.Range("I4").Value = 43.333332
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 43.333332
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 71.666668
$ws.Range("N4").Value = -330

$ws.Range("H20").Value = 1316.5
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 1633
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 1633
$ws.Range("M20").Value = -753
$ws.Range("N20").Value = -2127

$ws.Range("H22").Value = 485.5
$ws.Range("I22").Value = 485.5
$ws.Range("K22").Value = 485.5
$ws.Range("M22").Value = -312.5

$ws.Range("H99").Value = 2615.4285
$ws.Range("I99").Value = 1952.25
$ws.Range("K99").Value = 1952.25
$ws.Range("M99").Value = -454.25

$ws.Range("H134").Value = 5222.8184
$ws.Range("I134").Value = 5245.1
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 15735.3
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -13200.3
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2966.3333
$ws.Range("J62").Value = 2949.5
$ws.Range("L62").Value = 2949.5
$ws.Range("N62").Value = -4197.5

$ws.Range("H65").Value = 2966.3333
$ws.Range("J65").Value = 2949.5
$ws.Range("L65").Value = 14747.5
$ws.Range("N65").Value = -20987.5

$ws.Range("H99").Value = 2871.7856
$ws.Range("J99").Value = 3665
$ws.Range("L99").Value = 3665
$ws.Range("N99").Value = -6661

$ws.Range("H126").Value = 2871.7856
$ws.Range("J126").Value = 3665
$ws.Range("L126").Value = 10995
$ws.Range("N126").Value = -15935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 200000260
$ws.Range("I4").Value = 254.66667
$ws.Range("K4").Value = 764.00001
$ws.Range("M4").Value = -652.00001

$ws.Range("H12").Value = 157942.62
$ws.Range("I12").Value = 600022.2
$ws.Range("J12").Value = 57.07143
$ws.Range("K12").Value = 1800066.6
$ws.Range("L12").Value = 171.21429
$ws.Range("M12").Value = -1799893.6
$ws.Range("N12").Value = -517.21429

$ws.Range("H60").Value = 3976.8
$ws.Range("I60").Value = 2971.25
$ws.Range("K60").Value = 8913.75
$ws.Range("M60").Value = -8662.75

$ws.Range("H86").Value = 979.4167
$ws.Range("I86").Value = 1071.6666
$ws.Range("J86").Value = 887.1667
$ws.Range("K86").Value = 3214.9998
$ws.Range("L86").Value = 2661.5001
$ws.Range("M86").Value = -2028.9998
$ws.Range("N86").Value = -5033.5001

$ws.Range("H89").Value = 979.4167
$ws.Range("I89").Value = 1071.6666
$ws.Range("J89").Value = 887.1667
$ws.Range("K89").Value = 9644.999400000001
$ws.Range("L89").Value = 7984.5003
$ws.Range("M89").Value = -3716.999400000001
$ws.Range("N89").Value = -19840.5003

$ws.Range("H103").Value = 23118.818
$ws.Range("I103").Value = 41913.332
$ws.Range("J103").Value = 565.4
$ws.Range("K103").Value = 125739.996
$ws.Range("L103").Value = 1696.2
$ws.Range("M103").Value = -124860.996
$ws.Range("N103").Value = -3454.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8376
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 10166.333
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 10166.333
$ws.Range("M80").Value = -2007
$ws.Range("N80").Value = -12162.333

$ws.Range("H83").Value = 8376
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 10166.333
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 50831.665
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = -60815.665

$ws.Range("H102").Value = 1002.875
$ws.Range("I102").Value = 1002.875
$ws.Range("K102").Value = 1002.875
$ws.Range("M102").Value = 619.125

$ws.Range("H132").Value = 3640.6897
$ws.Range("I132").Value = 3556.4285
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 10669.2855
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -8139.2855
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 3766.611
$ws.Range("J46").Value = 4249.875
$ws.Range("L46").Value = 4249.875
$ws.Range("N46").Value = -4625.875

$ws.Range("H100").Value = 1433
$ws.Range("I100").Value = 1433
$ws.Range("K100").Value = 1433
$ws.Range("M100").Value = -892

$ws.Range("H123").Value = 25000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 25000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800
$ws.Range("M123").ClearContents()

$ws.Range("H132").Value = 7352.6665
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H135").Value = 91497.5
$ws.Range("J135").Value = 91497.5
$ws.Range("L135").Value = 91497.5
$ws.Range("N135").Value = -101637.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 20000
$ws.Range("K61").Value = 20000
$ws.Range("M61").Value = -19708

$ws.Range("H107").Value = 1595.2222
$ws.Range("I107").Value = 1651.6
$ws.Range("K107").Value = 4954.799999999999
$ws.Range("M107").Value = -3034.799999999999

$ws.Range("H122").Value = 2683.1667
$ws.Range("I122").Value = 2518
$ws.Range("K122").Value = 7554
$ws.Range("M122").Value = -5104

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
